$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 464, shifting existing rows (464-565) down to (465-566)
$ws.Rows.Item(464).Insert()

# Populate the newly inserted row 464 with the new record.
$ws.Cells.Item(464, 1).Value = 10
$ws.Cells.Item(464, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(464, 3).Value = "La Araucanía"
$ws.Cells.Item(464, 4).Value = 45173
$ws.Cells.Item(464, 5).Value = 9
$ws.Cells.Item(464, 6).Value = 100112009
$ws.Cells.Item(464, 7).Value = "Acelga"
$ws.Cells.Item(464, 8).Value = "Sin especificar"
$ws.Cells.Item(464, 9).Value = "Primera"
$ws.Cells.Item(464, 10).Value = 100
$ws.Cells.Item(464, 11).Value = 8000
$ws.Cells.Item(464, 12).Value = 8000
$ws.Cells.Item(464, 13).Value = 8000
$ws.Cells.Item(464, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(464, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(464, 16).Value = 667
$ws.Cells.Item(464, 17).Value = 12
$ws.Cells.Item(464, 18).Value = 'Hortaliza'
